$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "year" column header to "season_ending_year" (column B)
$ws.Range("B1").Value = "season_ending_year"

# Add the new "calendar_year" column header in AM1, matching the existing
# header formatting (bold, centered, bordered) by copying an adjacent
# header cell's format first.
$ws.Range("AL1").Copy() | Out-Null
$ws.Range("AM1").PasteSpecial(-4122) | Out-Null
$ws.Range("AM1").Value = "calendar_year"

# Find the last data row (row 1 is the header row)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row()

# Fill birth_year (column E) per player, and calendar_year (column AM) as
# a copy of the season_ending_year value (column B) for every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $player = $ws.Cells.Item($r, 4).Value()
    if ($player -eq "LeBron James") {
        $ws.Cells.Item($r, 5).Value = 1985
    } elseif ($player -eq "Michael Jordan") {
        $ws.Cells.Item($r, 5).Value = 1964
    }

    $yearValue = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($r, 39).Value = [int]$yearValue
}
